$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -7.884599999999999
$ws.Range("D7").Value = -7.640299999999992
$ws.Range("C8").Value = -11.02549999999999
$ws.Range("B12").Value = 5.547499999999997
$ws.Range("C12").Value = -14.77440000000002
$ws.Range("C14").Value = -12.3454
$ws.Range("D19").Value = -8.432299999999993
$ws.Range("D21").Value = -7.702699999999997
$ws.Range("C22").Value = -10.43609999999999
$ws.Range("D24").Value = -8.438499999999992
